$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H2").Value = 892.4286
$ws.Range("J2").Value = 1169.25
$ws.Range("L2").Value = 1169.25
$ws.Range("N2").Value = -1395.25

$ws.Range("H12").Value = 6381.0625
$ws.Range("I12").Value = 6800.8
$ws.Range("J12").Value = 85
$ws.Range("K12").Value = 6800.8
$ws.Range("L12").Value = 85
$ws.Range("M12").Value = -6630.8
$ws.Range("N12").Value = -425

$ws.Range("H17").Value = 999999
$ws.Range("I17").Value = 0
$ws.Range("J17").Value = 999999
$ws.Range("K17").Value = 0
$ws.Range("L17").Value = 2999997
$ws.Range("M17").ClearContents()
$ws.Range("N17").Value = -3000333

$ws.Range("H19").Value = 547.7917
$ws.Range("I19").Value = 391.25
$ws.Range("J19").Value = 626.0625
$ws.Range("K19").Value = 391.25
$ws.Range("L19").Value = 626.0625
$ws.Range("M19").Value = -216.25
$ws.Range("N19").Value = -976.0625

$ws.Range("H106").Value = 2578.9048
$ws.Range("I106").Value = 2482
$ws.Range("K106").Value = 2482
$ws.Range("M106").Value = -1851

$ws.Range("H137").Value = 2677.3333
$ws.Range("I137").Value = 1442.4286
$ws.Range("K137").Value = 4327.2858
$ws.Range("M137").Value = -1777.2858

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 465487.62
$ws.Range("I2").Value = 1042258.5
$ws.Range("J2").Value = 4070.95
$ws.Range("K2").Value = 1042258.5
$ws.Range("L2").Value = 4070.95
$ws.Range("M2").Value = -1042145.5
$ws.Range("N2").Value = -4296.95

$ws.Range("H44").Value = 53950
$ws.Range("J44").Value = 53950
$ws.Range("L44").Value = 53950
$ws.Range("N44").Value = -54926

$ws.Range("H45").Value = 1611.6364
$ws.Range("I45").Value = 1542.9
$ws.Range("J45").Value = 2299
$ws.Range("K45").Value = 1542.9
$ws.Range("L45").Value = 2299
$ws.Range("M45").Value = -1165.9
$ws.Range("N45").Value = -3053

$ws.Range("H61").Value = 37041600
$ws.Range("I61").Value = 41671310
$ws.Range("J61").Value = 3900
$ws.Range("K61").Value = 41671310
$ws.Range("L61").Value = 3900
$ws.Range("M61").Value = -41671098
$ws.Range("N61").Value = -4324

$ws.Range("H116").Value = 465487.62
$ws.Range("I116").Value = 1042258.5
$ws.Range("J116").Value = 4070.95
$ws.Range("K116").Value = 1042258.5
$ws.Range("L116").Value = 4070.95
$ws.Range("M116").Value = -1039964.5
$ws.Range("N116").Value = -8658.950000000001

$ws.Range("H132").Value = 1854609
$ws.Range("I132").Value = 1963676.6
$ws.Range("J132").Value = 459.33334
$ws.Range("K132").Value = 5891029.800000001
$ws.Range("L132").Value = 1378.00002
$ws.Range("M132").Value = -5888499.800000001
$ws.Range("N132").Value = -6438.000019999999

$ws.Range("H136").Value = 37041600
$ws.Range("I136").Value = 41671310
$ws.Range("J136").Value = 3900
$ws.Range("K136").Value = 125013930
$ws.Range("L136").Value = 11700
$ws.Range("M136").Value = -125011380
$ws.Range("N136").Value = -16800

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 2553.4546
$ws.Range("I20").Value = 2267.2942
$ws.Range("J20").Value = 3526.4
$ws.Range("K20").Value = 2267.2942
$ws.Range("L20").Value = 3526.4
$ws.Range("M20").Value = -2020.2942
$ws.Range("N20").Value = -4020.4

$ws.Range("H86").Value = 5003
$ws.Range("I86").Value = 5003
$ws.Range("K86").Value = 5003
$ws.Range("M86").Value = -3880

$ws.Range("H89").Value = 5003
$ws.Range("I89").Value = 5003
$ws.Range("K89").Value = 25015
$ws.Range("M89").Value = -19399

$ws.Range("H107").Value = 47257.953
$ws.Range("I107").Value = 1333.1428
$ws.Range("K107").Value = 1333.1428
$ws.Range("M107").Value = 586.8571999999999

$ws.Range("H134").Value = 16132170
$ws.Range("I134").Value = 21741170
$ws.Range("J134").Value = 6293.5
$ws.Range("K134").Value = 65223510
$ws.Range("L134").Value = 18880.5
$ws.Range("M134").Value = -65220975
$ws.Range("N134").Value = -23950.5

$ws.Range("H138").Value = 122633
$ws.Range("J138").Value = 122633
$ws.Range("L138").Value = 122633
$ws.Range("N138").Value = -132913

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H22").Value = 8861.833000000001
$ws.Range("J22").Value = 611.5
$ws.Range("L22").Value = 611.5
$ws.Range("N22").Value = -1311.5

$ws.Range("H47").Value = 35442
$ws.Range("I47").Value = 30885
$ws.Range("K47").Value = 30885
$ws.Range("M47").Value = -30319

$ws.Range("H105").Value = 2501174.5
$ws.Range("I105").Value = 3334374.5
$ws.Range("J105").Value = 1575
$ws.Range("K105").Value = 3334374.5
$ws.Range("L105").Value = 1575
$ws.Range("M105").Value = -3332627.5
$ws.Range("N105").Value = -5069

$ws.Range("H132").Value = 27027948
$ws.Range("I132").Value = 31250830
$ws.Range("J132").Value = 1509.6
$ws.Range("K132").Value = 93752490
$ws.Range("L132").Value = 4528.799999999999
$ws.Range("M132").Value = -93749960
$ws.Range("N132").Value = -9588.799999999999

$ws.Range("H134").Value = 20834168
$ws.Range("I134").Value = 21739976
$ws.Range("K134").Value = 65219928
$ws.Range("M134").Value = -65217393

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 316.92307
$ws.Range("I2").Value = 243.71428
$ws.Range("J2").Value = 402.33334
$ws.Range("K2").Value = 1462.28568
$ws.Range("L2").Value = 2414.00004
$ws.Range("M2").Value = -1349.28568
$ws.Range("N2").Value = -2640.00004

$ws.Range("H7").Value = 770056.9
$ws.Range("I7").Value = 1250760.9
$ws.Range("J7").Value = 930.4
$ws.Range("K7").Value = 3752282.7
$ws.Range("L7").Value = 2791.2
$ws.Range("M7").Value = -3752170.7
$ws.Range("N7").Value = -3015.2

$ws.Range("H14").Value = 460.875
$ws.Range("I14").Value = 460.875
$ws.Range("K14").Value = 1382.625
$ws.Range("M14").Value = -1209.625

$ws.Range("H23").Value = 115.6
$ws.Range("I23").Value = 400
$ws.Range("J23").Value = 44.5
$ws.Range("K23").Value = 1200
$ws.Range("L23").Value = 133.5
$ws.Range("M23").Value = -965
$ws.Range("N23").Value = -603.5

$ws.Range("H62").Value = 5299.6665
$ws.Range("J62").Value = 6699.5
$ws.Range("L62").Value = 20098.5
$ws.Range("N62").Value = -21470.5

$ws.Range("H65").Value = 5299.6665
$ws.Range("J65").Value = 6699.5
$ws.Range("L65").Value = 60295.5
$ws.Range("N65").Value = -67159.5

$ws.Range("H107").Value = 1268.7037
$ws.Range("I107").Value = 348.1
$ws.Range("J107").Value = 1810.2354
$ws.Range("K107").Value = 1044.3
$ws.Range("L107").Value = 5430.706200000001
$ws.Range("M107").Value = 875.6999999999998
$ws.Range("N107").Value = -9270.706200000001

$ws.Range("H131").Value = 1914.25
$ws.Range("J131").Value = 1997.2941
$ws.Range("L131").Value = 5991.8823
$ws.Range("N131").Value = -16071.8823

$ws.Range("H132").Value = 1609.381
$ws.Range("I132").Value = 1665.7858
$ws.Range("J132").Value = 1496.5714
$ws.Range("K132").Value = 14992.0722
$ws.Range("L132").Value = 13469.1426
$ws.Range("M132").Value = -12462.0722
$ws.Range("N132").Value = -18529.1426

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 59162.555
$ws.Range("I113").Value = 85755
$ws.Range("K113").Value = 85755
$ws.Range("M113").Value = -83585

$ws.Range("H126").Value = 4502.7744
$ws.Range("I126").Value = 4536.2
$ws.Range("K126").Value = 13608.6
$ws.Range("M126").Value = -11138.6

$ws.Range("H132").Value = 4465281.5
$ws.Range("I132").Value = 4808564.5
$ws.Range("K132").Value = 14425693.5
$ws.Range("M132").Value = -14423163.5

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5120.8
$ws.Range("I7").Value = 4558.0713
$ws.Range("K7").Value = 4558.0713
$ws.Range("M7").Value = -4446.0713

$ws.Range("H16").Value = 3073.8696
$ws.Range("I16").Value = 791.9
$ws.Range("K16").Value = 791.9
$ws.Range("M16").Value = -621.9

$ws.Range("H46").Value = 757.6
$ws.Range("J46").Value = 784.2857
$ws.Range("L46").Value = 784.2857
$ws.Range("N46").Value = -1160.2857

$ws.Range("H55").Value = 331.16327
$ws.Range("I55").Value = 131.17392
$ws.Range("K55").Value = 131.17392
$ws.Range("M55").Value = 41.82607999999999

$ws.Range("H126").Value = 5120.8
$ws.Range("I126").Value = 4558.0713
$ws.Range("K126").Value = 13674.2139
$ws.Range("M126").Value = -11204.2139

$ws.Range("H136").Value = 2867.9375
$ws.Range("I136").Value = 1633.3334
$ws.Range("K136").Value = 4900.0002
$ws.Range("M136").Value = -2350.0002

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 3499.5
$ws.Range("I2").Value = 3499.5
$ws.Range("K2").Value = 3499.5
$ws.Range("M2").Value = -3387.5

$ws.Range("H132").Value = 10207649
$ws.Range("I132").Value = 12198897
$ws.Range("J132").Value = 2498.875
$ws.Range("K132").Value = 36596691
$ws.Range("L132").Value = 7496.625
$ws.Range("M132").Value = -36594161
$ws.Range("N132").Value = -12556.625

$ws.Range("H136").Value = 23811320
$ws.Range("J136").Value = 4977.4
$ws.Range("L136").Value = 14932.2
$ws.Range("N136").Value = -20032.2
